$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting old B->C and C->D
$ws.Columns("B:B").Insert()

# New column B width should match column A's width (75.81640625)
$aWidth = $ws.Columns("A:A").ColumnWidth
$ws.Columns("B:B").ColumnWidth = $aWidth

# Header for new column
$ws.Range("B1").Value2 = "StatQuery"

# New query text for row 2, column B, with same wrap-text style as A2
$ws.Range("B2").Value2 = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Adenocarcinoma of the colon'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

$ws.Range("B2").WrapText = $true

# Update the selection to A2
$ws.Range("A2").Select()
